$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("BPaFF-BITPTaP")
$ws2.Range("B15").Formula = "=B11"

$ws3 = $wb.Worksheets.Item("BPaFF-BDTPTPF")
$ws3.Range("B15").Formula = "=B11"
